# paises.xlsx refresh (27 Mar 2020, 00:42 -> 01:12)
#
# The sheet "Pais" lists one country per row (column A) with its stats in
# columns B:H (Casos totales, Nuevos casos, Casos activos, Recuperados,
# Casos criticos, Muertes hoy, Muertes). Column A cells are shared-string
# references, so a handful of rows whose countries swapped position in the
# refreshed export effectively show a different country name even though
# their row number did not change. Each such row also got new stats from
# the refreshed data pull; a few other rows only had their stats refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "last updated" banner (row 1) ---------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 01:12"

# --- Totals row (Estados Unidos, row 4) refreshed ------------------------
$ws.Range("B4").Value = 84911
$ws.Range("C4").Value = 16700
$ws.Range("E4").Value = 81757
$ws.Range("G4").Value = 263
$ws.Range("H4").Value = 1290

# --- Australia (row 21) refreshed ----------------------------------------
$ws.Range("B21").Value = 3050
$ws.Range("C21").Value = 374
$ws.Range("E21").Value = 2867

# --- Rows 64-68: countries reshuffled (Argelia..Nueva Zelanda block) -----
# Nueva Zelanda moves up in front of Argelia; every row here gets the
# data that used to belong to the row above it, and row 64 gets fresh data.
$ws.Range("A64").Value = "Nueva Zelanda"
$ws.Range("B64").Value = 368
$ws.Range("C64").Value = 85
$ws.Range("D64").Value = 37
$ws.Range("E64").Value = 331
$ws.Range("F64").Value = 1
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0

$ws.Range("A65").Value = "Argelia"
$ws.Range("B65").Value = 367
$ws.Range("C65").Value = 65
$ws.Range("D65").Value = 29
$ws.Range("E65").Value = 313
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 25

$ws.Range("A66").Value = "Emiratos Arabes Unidos"
$ws.Range("B66").Value = 333
$ws.Range("C66").Value = 0
$ws.Range("D66").Value = 52
$ws.Range("E66").Value = 279
$ws.Range("F66").Value = 2
$ws.Range("H66").Value = 2

$ws.Range("A67").Value = "Lituania"
$ws.Range("B67").Value = 299
$ws.Range("D67").Value = 1
$ws.Range("E67").Value = 294
$ws.Range("F67").Value = 1
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 4

$ws.Range("A68").Value = "Armenia"
$ws.Range("B68").Value = 290
$ws.Range("C68").Value = 25
$ws.Range("D68").Value = 18
$ws.Range("E68").Value = 271
$ws.Range("F68").Value = 6
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 1

# --- Rows 137-138: Togo moves in front of Madagascar ----------------------
$ws.Range("A137").Value = "Togo"
$ws.Range("B137").Value = 24
$ws.Range("C137").Value = 1
$ws.Range("D137").Value = 1

$ws.Range("A138").Value = "Madagascar"
$ws.Range("C138").Value = 4
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 23

# --- Rows 145-146: Tanzania moves in front of El Salvador ------------------
$ws.Range("A145").Value = "Tanzania"
$ws.Range("C145").Value = 0

$ws.Range("A146").Value = "El Salvador"
$ws.Range("C146").Value = 4

# --- Rows 148-153: Guinea Ecuatorial / Etiopia / San Martin / Yibuti /
#     Dominica / Mongolia reshuffle -----------------------------------------
$ws.Range("A148").Value = "Guinea Ecuatorial"
$ws.Range("C148").Value = 3

$ws.Range("A149").Value = "Etiopia"
$ws.Range("C149").Value = 0

$ws.Range("A150").Value = "San Martin (Parte Francesa)"
$ws.Range("C150").Value = 0

$ws.Range("A151").Value = "Republica de Yibuti"

$ws.Range("A152").Value = "Dominica"

$ws.Range("A153").Value = "Mongolia"
$ws.Range("C153").Value = 1

# --- Rows 160-163: Granada / Seychelles move in front of Antigua y
#     Barbuda / Mozambique --------------------------------------------------
$ws.Range("A160").Value = "Granada"
$ws.Range("C160").Value = 6

$ws.Range("A161").Value = "Seychelles"
$ws.Range("C161").Value = 0

$ws.Range("A162").Value = "Antigua y Barbuda"
$ws.Range("C162").Value = 4

$ws.Range("A163").Value = "Mozambique"
$ws.Range("C163").Value = 2
